$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 5.256564333333333
$ws.Cells.Item(2, 8).Value = 15.769693
$ws.Cells.Item(2, 9).Value = 0.003747859920520347
$ws.Cells.Item(2, 10).Value = 0.003747859920520347
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 122.328922
$ws.Cells.Item(2, 14).Value = 366.986766
$ws.Cells.Item(2, 15).Value = 0.9783373008518612
$ws.Cells.Item(2, 16).Value = 0.9783373008518613
$ws.Cells.Item(2, 17).Value = 643.0298483203153
$ws.Cells.Item(2, 18).Value = 5787.268634882838
$ws.Cells.Item(2, 19).Value = 0.003666671158612747
$ws.Cells.Item(2, 20).Value = 0.003666671158612747
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 5.256564333333333
$ws.Cells.Item(3, 8).Value = 15.769693
$ws.Cells.Item(3, 9).Value = 0.003747859920520347
$ws.Cells.Item(3, 10).Value = 0.003747859920520347
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.3863573333333334
$ws.Cells.Item(3, 14).Value = 1.159072
$ws.Cells.Item(3, 15).Value = 0.003089929874945324
$ws.Cells.Item(3, 16).Value = 0.003089929874945324
$ws.Cells.Item(3, 17).Value = 2.030912178321778
$ws.Cells.Item(3, 18).Value = 18.278209604896
$ws.Cells.Item(3, 19).Value = [double]"1.158062433552603E-05"
$ws.Cells.Item(3, 20).Value = [double]"1.158062433552602E-05"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 5.256564333333333
$ws.Cells.Item(4, 8).Value = 15.769693
$ws.Cells.Item(4, 9).Value = 0.003747859920520347
$ws.Cells.Item(4, 10).Value = 0.003747859920520347
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 2.322294
$ws.Cells.Item(4, 14).Value = 6.966882000000001
$ws.Cells.Item(4, 15).Value = 0.0185727692731934
$ws.Cells.Item(4, 16).Value = 0.0185727692731934
$ws.Cells.Item(4, 17).Value = 12.207287811914
$ws.Cells.Item(4, 18).Value = 109.865590307226
$ws.Cells.Item(4, 19).Value = [double]"6.960813757207338E-05"
$ws.Cells.Item(4, 20).Value = [double]"6.960813757207338E-05"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 1312.703450666667
$ws.Cells.Item(5, 8).Value = 3938.110352
$ws.Cells.Item(5, 9).Value = 0.93593996730609
$ws.Cells.Item(5, 10).Value = 0.9359399673060897
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 122.328922
$ws.Cells.Item(5, 14).Value = 366.986766
$ws.Cells.Item(5, 15).Value = 0.9783373008518612
$ws.Cells.Item(5, 16).Value = 0.9783373008518613
$ws.Cells.Item(5, 17).Value = 160581.5980257335
$ws.Cells.Item(5, 18).Value = 1445234.382231602
$ws.Cells.Item(5, 19).Value = 0.9156649813736193
$ws.Cells.Item(5, 20).Value = 0.9156649813736192
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1312.703450666667
$ws.Cells.Item(6, 8).Value = 3938.110352
$ws.Cells.Item(6, 9).Value = 0.93593996730609
$ws.Cells.Item(6, 10).Value = 0.9359399673060897
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.3863573333333334
$ws.Cells.Item(6, 14).Value = 1.159072
$ws.Cells.Item(6, 15).Value = 0.003089929874945324
$ws.Cells.Item(6, 16).Value = 0.003089929874945324
$ws.Cells.Item(6, 17).Value = 507.1726046570383
$ws.Cells.Item(6, 18).Value = 4564.553441913344
$ws.Cells.Item(6, 19).Value = 0.002891988866134437
$ws.Cells.Item(6, 20).Value = 0.002891988866134436
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1312.703450666667
$ws.Cells.Item(7, 8).Value = 3938.110352
$ws.Cells.Item(7, 9).Value = 0.93593996730609
$ws.Cells.Item(7, 10).Value = 0.9359399673060897
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 2.322294
$ws.Cells.Item(7, 14).Value = 6.966882000000001
$ws.Cells.Item(7, 15).Value = 0.0185727692731934
$ws.Cells.Item(7, 16).Value = 0.0185727692731934
$ws.Cells.Item(7, 17).Value = 3048.483347262496
$ws.Cells.Item(7, 18).Value = 27436.35012536247
$ws.Cells.Item(7, 19).Value = 0.01738299706633619
$ws.Cells.Item(7, 20).Value = 0.01738299706633618
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 84.59089266666666
$ws.Cells.Item(8, 8).Value = 253.772678
$ws.Cells.Item(8, 9).Value = 0.06031217277338979
$ws.Cells.Item(8, 10).Value = 0.06031217277338978
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 122.328922
$ws.Cells.Item(8, 14).Value = 366.986766
$ws.Cells.Item(8, 15).Value = 0.9783373008518612
$ws.Cells.Item(8, 16).Value = 0.9783373008518613
$ws.Cells.Item(8, 17).Value = 10347.91271093104
$ws.Cells.Item(8, 18).Value = 93131.21439837934
$ws.Cells.Item(8, 19).Value = 0.05900564831962928
$ws.Cells.Item(8, 20).Value = 0.05900564831962928
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 84.59089266666666
$ws.Cells.Item(9, 8).Value = 253.772678
$ws.Cells.Item(9, 9).Value = 0.06031217277338979
$ws.Cells.Item(9, 10).Value = 0.06031217277338978
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.3863573333333334
$ws.Cells.Item(9, 14).Value = 1.159072
$ws.Cells.Item(9, 15).Value = 0.003089929874945324
$ws.Cells.Item(9, 16).Value = 0.003089929874945324
$ws.Cells.Item(9, 17).Value = 32.68231171497956
$ws.Cells.Item(9, 18).Value = 294.140805434816
$ws.Cells.Item(9, 19).Value = 0.0001863603844753611
$ws.Cells.Item(9, 20).Value = 0.000186360384475361
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 84.59089266666666
$ws.Cells.Item(10, 8).Value = 253.772678
$ws.Cells.Item(10, 9).Value = 0.06031217277338979
$ws.Cells.Item(10, 10).Value = 0.06031217277338978
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 2.322294
$ws.Cells.Item(10, 14).Value = 6.966882000000001
$ws.Cells.Item(10, 15).Value = 0.0185727692731934
$ws.Cells.Item(10, 16).Value = 0.0185727692731934
$ws.Cells.Item(10, 17).Value = 196.444922494444
$ws.Cells.Item(10, 18).Value = 1768.004302449996
$ws.Cells.Item(10, 19).Value = 0.001120164069285146
$ws.Cells.Item(10, 20).Value = 0.001120164069285146
